{"js": "const pairs = [\n  [\"2023-03-15 Wednesday\", \"2023-03-16 Thursday\"],\n  [\"69-20=49\", \"26+0=26\"],\n  [\"26+35=61\", \"23+52=75\"],\n  [\"9+30=39\", \"30+23=53\"],\n  [\"93-3=90\", \"13+28=41\"],\n  [\"54-38=16\", \"79-27=52\"],\n  [\"59+9=68\", \"71+4=75\"],\n  [\"82-52=30\", \"18+10=28\"],\n  [\"91-48=43\", \"89-49=40\"],\n  [\"26+14=40\", \"44-1=43\"],\n  [\"16+31=47\", \"69-8=61\"],\n  [\"19+9=28\", \"38+17=55\"],\n  [\"95-58=37\", \"28+24=52\"],\n  [\"31-14=17\", \"71+26=97\"],\n  [\"42+44=86\", \"17+31=48\"],\n  [\"65+13=78\", \"33+30=63\"],\n  [\"99-57=42\", \"28+5=33\"],\n  [\"14-9=5\", \"73-3=70\"],\n  [\"28+6=34\", \"49+3=52\"],\n  [\"71-14=57\", \"52+22=74\"],\n  [\"97-79=18\", \"64+35=99\"],\n  [\"30+30=60\", \"76+4=80\"],\n  [\"55+41=96\", \"17+26=43\"],\n  [\"39+12=51\", \"83-32=51\"],\n  [\"65+29=94\", \"27-2=25\"],\n  [\"80-65=15\", \"66+17=83\"],\n  [\"55-36=19\", \"18+53=71\"],\n  [\"7+67=74\", \"22+26=48\"],\n  [\"0+38=38\", \"58+33=91\"],\n  [\"60-6=54\", \"16-9=7\"],\n  [\"46+39=85\", \"15+6=21\"],\n  [\"91+2=93\", \"60-13=47\"],\n  [\"61-12=49\", \"11+20=31\"],\n  [\"52+8=60\", \"90-68=22\"],\n  [\"18+48=66\", \"71-60=11\"],\n  [\"31+37=68\", \"95-83=12\"],\n  [\"64+25=89\", \"96-55=41\"],\n  [\"92-53=39\", \"9+17=26\"],\n  [\"64-9=55\", \"21-9=12\"],\n  [\"34+2=36\", \"28+25=53\"],\n  [\"98-61=37\", \"34+31=65\"],\n  [\"20+10=30\", \"88+7=95\"],\n  [\"63-26=37\", \"41+52=93\"],\n  [\"67-35=32\", \"27+36=63\"],\n  [\"11-4=7\", \"62+16=78\"],\n  [\"95-51=44\", \"21+32=53\"],\n  [\"36-12=24\", \"70-42=28\"],\n  [\"11+38=49\", \"78+5=83\"],\n  [\"84-16=68\", \"87-36=51\"],\n  [\"52-46=6\", \"66-27=39\"],\n  [\"73-19=54\", \"25+60=85\"],\n  [\"51-28=23\", \"59+8=67\"],\n  [\"75-11=64\", \"19+24=43\"],\n  [\"20+8=28\", \"24+16=40\"],\n  [\"75-2=73\", \"65-59=6\"],\n  [\"42-9=33\", \"45+48=93\"],\n  [\"19+34=53\", \"74-9=65\"],\n  [\"14+41=55\", \"65+3=68\"],\n  [\"50-26=24\", \"41+2=43\"],\n  [\"26+33=59\", \"79-72=7\"],\n  [\"65+34=99\", \"2+89=91\"],\n  [\"91-2=89\", \"16+25=41\"],\n  [\"28+68=96\", \"68-60=8\"],\n  [\"98-7=91\", \"90-4=86\"],\n  [\"21+48=69\", \"85-57=28\"],\n  [\"56+23=79\", \"95-7=88\"],\n  [\"32+42=74\", \"79-25=54\"],\n  [\"81-53=28\", \"8+64=72\"],\n  [\"23+22=45\", \"9+73=82\"],\n  [\"34+16=50\", \"48-44=4\"],\n  [\"16+56=72\", \"9+54=63\"],\n  [\"3+88=91\", \"48+26=74\"],\n  [\"39-5=34\", \"49-31=18\"],\n  [\"9+9=18\", \"87-35=52\"],\n  [\"37+2=39\", \"31+24=55\"],\n  [\"25+55=80\", \"77+15=92\"],\n  [\"25+30=55\", \"95-66=29\"],\n  [\"56+26=82\", \"63+31=94\"],\n  [\"78-17=61\", \"84-44=40\"],\n  [\"83-12=71\", \"80-6=74\"],\n  [\"22+32=54\", \"87-21=66\"],\n  [\"69+23=92\", \"54-33=21\"],\n  [\"23-7=16\", \"80-54=26\"],\n  [\"60-30=30\", \"4+35=39\"],\n  [\"18+74=92\", \"57+33=90\"],\n  [\"92-51=41\", \"89-86=3\"],\n  [\"72-48=24\", \"16+53=69\"],\n  [\"25+28=53\", \"62-54=8\"],\n  [\"91-10=81\", \"30+7=37\"],\n  [\"67+13=80\", \"81-35=46\"],\n  [\"86-33=53\", \"77-41=36\"],\n  [\"12+49=61\", \"50+25=75\"],\n  [\"2+10=12\", \"1+5=6\"],\n  [\"8+44=52\", \"41+39=80\"],\n  [\"44+24=68\", \"87-18=69\"],\n  [\"10-4=6\", \"10+42=52\"],\n  [\"46+53=99\", \"20+73=93\"],\n  [\"40+42=82\", \"3+28=31\"],\n  [\"33-13=20\", \"90+6=96\"],\n  [\"16+55=71\", \"96-79=17\"],\n  [\"37-3=34\", \"86-1=85\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n    @('2023-03-15 Wednesday', '2023-03-16 Thursday'),\n    @('69-20=49', '26+0=26'),\n    @('26+35=61', '23+52=75'),\n    @('9+30=39', '30+23=53'),\n    @('93-3=90', '13+28=41'),\n    @('54-38=16', '79-27=52'),\n    @('59+9=68', '71+4=75'),\n    @('82-52=30', '18+10=28'),\n    @('91-48=43', '89-49=40'),\n    @('26+14=40', '44-1=43'),\n    @('16+31=47', '69-8=61'),\n    @('19+9=28', '38+17=55'),\n    @('95-58=37', '28+24=52'),\n    @('31-14=17', '71+26=97'),\n    @('42+44=86', '17+31=48'),\n    @('65+13=78', '33+30=63'),\n    @('99-57=42', '28+5=33'),\n    @('14-9=5', '73-3=70'),\n    @('28+6=34', '49+3=52'),\n    @('71-14=57', '52+22=74'),\n    @('97-79=18', '64+35=99'),\n    @('30+30=60', '76+4=80'),\n    @('55+41=96', '17+26=43'),\n    @('39+12=51', '83-32=51'),\n    @('65+29=94', '27-2=25'),\n    @('80-65=15', '66+17=83'),\n    @('55-36=19', '18+53=71'),\n    @('7+67=74', '22+26=48'),\n    @('0+38=38', '58+33=91'),\n    @('60-6=54', '16-9=7'),\n    @('46+39=85', '15+6=21'),\n    @('91+2=93', '60-13=47'),\n    @('61-12=49', '11+20=31'),\n    @('52+8=60', '90-68=22'),\n    @('18+48=66', '71-60=11'),\n    @('31+37=68', '95-83=12'),\n    @('64+25=89', '96-55=41'),\n    @('92-53=39', '9+17=26'),\n    @('64-9=55', '21-9=12'),\n    @('34+2=36', '28+25=53'),\n    @('98-61=37', '34+31=65'),\n    @('20+10=30', '88+7=95'),\n    @('63-26=37', '41+52=93'),\n    @('67-35=32', '27+36=63'),\n    @('11-4=7', '62+16=78'),\n    @('95-51=44', '21+32=53'),\n    @('36-12=24', '70-42=28'),\n    @('11+38=49', '78+5=83'),\n    @('84-16=68', '87-36=51'),\n    @('52-46=6', '66-27=39'),\n    @('73-19=54', '25+60=85'),\n    @('51-28=23', '59+8=67'),\n    @('75-11=64', '19+24=43'),\n    @('20+8=28', '24+16=40'),\n    @('75-2=73', '65-59=6'),\n    @('42-9=33', '45+48=93'),\n    @('19+34=53', '74-9=65'),\n    @('14+41=55', '65+3=68'),\n    @('50-26=24', '41+2=43'),\n    @('26+33=59', '79-72=7'),\n    @('65+34=99', '2+89=91'),\n    @('91-2=89', '16+25=41'),\n    @('28+68=96', '68-60=8'),\n    @('98-7=91', '90-4=86'),\n    @('21+48=69', '85-57=28'),\n    @('56+23=79', '95-7=88'),\n    @('32+42=74', '79-25=54'),\n    @('81-53=28', '8+64=72'),\n    @('23+22=45', '9+73=82'),\n    @('34+16=50', '48-44=4'),\n    @('16+56=72', '9+54=63'),\n    @('3+88=91', '48+26=74'),\n    @('39-5=34', '49-31=18'),\n    @('9+9=18', '87-35=52'),\n    @('37+2=39', '31+24=55'),\n    @('25+55=80', '77+15=92'),\n    @('25+30=55', '95-66=29'),\n    @('56+26=82', '63+31=94'),\n    @('78-17=61', '84-44=40'),\n    @('83-12=71', '80-6=74'),\n    @('22+32=54', '87-21=66'),\n    @('69+23=92', '54-33=21'),\n    @('23-7=16', '80-54=26'),\n    @('60-30=30', '4+35=39'),\n    @('18+74=92', '57+33=90'),\n    @('92-51=41', '89-86=3'),\n    @('72-48=24', '16+53=69'),\n    @('25+28=53', '62-54=8'),\n    @('91-10=81', '30+7=37'),\n    @('67+13=80', '81-35=46'),\n    @('86-33=53', '77-41=36'),\n    @('12+49=61', '50+25=75'),\n    @('2+10=12', '1+5=6'),\n    @('8+44=52', '41+39=80'),\n    @('44+24=68', '87-18=69'),\n    @('10-4=6', '10+42=52'),\n    @('46+53=99', '20+73=93'),\n    @('40+42=82', '3+28=31'),\n    @('33-13=20', '90+6=96'),\n    @('16+55=71', '96-79=17'),\n    @('37-3=34', '86-1=85'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n\nWrite-Output \"Replaced $($pairs.Count) items\""}
